$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.030.48"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'2.942.84"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'375.14"
$ws.Range("E5").Value = "  -0.65%  "
$ws.Range("D6").Value = "'102.28"
$ws.Range("E6").Value = "  -2.00%  "
$ws.Range("D7").Value = "'0.536"
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.583"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "'36.51"
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'0.0837"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("D13").Value = "'3.403.14"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("D14").Value = "'17.91"
$ws.Range("E14").Value = "  -2.15%  "
$ws.Range("D15").Value = "'7.36"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "'2.940.40"
$ws.Range("E16").Value = "  +0.21%  "
$ws.Range("D17").Value = "'0.981"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "'50.980.80"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("E19").Value = "  -5.35%  "
$ws.Range("D20").Value = "'7.18"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").Value = "'12.62"
$ws.Range("E21").Value = "  -1.64%  "
$ws.Range("D22").Value = "'0.0₃0957"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'264.25"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "'68.35"
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +2.94%  "
$ws.Range("D26").Value = "'8.36"
$ws.Range("E26").Value = "  +11.94%  "
$ws.Range("D27").Value = "'7.94"
$ws.Range("E27").Value = "  +9.73%  "
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("D30").Value = "'0.113"
$ws.Range("E30").Value = "  +1.27%  "
$ws.Range("D31").Value = "'25.63"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").Value = "'9.88"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "'50.83"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "'33.71"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("D35").Value = "'0.0449"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "'2.99"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.116"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "'16.45"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "'120.48"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("E44").Value = "  +3.64%  "
$ws.Range("E45").Value = "  -4.36%  "
$ws.Range("D46").Value = "'2.03"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").Value = "'3.26"
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D49").Value = "'1.979.66"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "'0.0342"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").Value = "'5.06"
$ws.Range("E51").Value = "  -0.33%  "
